# Update the "Execution Time" column (E) values on the Test Results sheet
# to reflect the report being regenerated at a later time, per the commit
# "Updated report to delete automatically".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "03/27/2025 06:24:01 PM"
$ws.Range("E3").Value = "03/27/2025 06:24:06 PM"
$ws.Range("E4").Value = "03/27/2025 06:24:07 PM"
$ws.Range("E5").Value = "03/27/2025 06:24:09 PM"
$ws.Range("E6").Value = "03/27/2025 06:24:09 PM"
$ws.Range("E7").Value = "03/27/2025 06:24:10 PM"
$ws.Range("E8").Value = "03/27/2025 06:24:10 PM"
